$d = $word.ActiveDocument
$quote = [char]8221

# ---------------------------------------------------------------------------
# Paragraph: "First fragment encounter" -> AH's dialogue line.
# Replace "know this artefact, but I have never studied this civilization in
# detail" with "can see the full artefact right in in front of me".
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("know this artefact, but I have never studied this civilization in detail", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find first fragment dialogue text to replace"
}
$rng1.Text = "can see the full artefact right in in front of me"
$newEnd1 = $rng1.End
# Toggling a character property and reverting it keeps this inserted phrase in
# its own run instead of it re-absorbing into the (identically formatted)
# text that precedes it.
$rng1.Bold = 1
$rng1.Bold = 0

# The ellipsis + closing curly quote that follow immediately are untouched,
# pre-existing runs; protect that boundary too so our edit doesn't silently
# swallow them into the freshly written run.
$tail1 = $d.Range($newEnd1, $d.Content.End)
$foundTail1 = $tail1.Find.Execute($quote, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundTail1) {
    $tail1.Bold = 1
    $tail1.Bold = 0
}

# ---------------------------------------------------------------------------
# Paragraph: AH's next line - "Ah, probably just my intuition. ..."
# Replace "check for details on the fragment”." with
# "analyse it in more detail.”" (note the closing quote now follows the
# period instead of preceding it).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("check for details on the fragment" + $quote + ".", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find second fragment dialogue text to replace"
}
$p2Start = $rng2.Paragraphs(1).Range.Start
$rng2.Text = "analyse it in more detail." + $quote
$newStart2 = $rng2.Start
$newEnd2 = $rng2.End
$rng2.Bold = 1
$rng2.Bold = 0

# Protect the boundary with the untouched 'AH: "' lead-in run at the start of
# the same paragraph so it doesn't get folded into our rewritten text.
$lead2 = $d.Range($p2Start, $newStart2)
if ($lead2.Start -lt $lead2.End) {
    $lead2.Bold = 1
    $lead2.Bold = 0
}

# Split the freshly inserted text into "analyse it in more detail." and the
# closing curly quote as two separate runs.
$closeQuote2 = $d.Range($newEnd2 - 1, $newEnd2)
$closeQuote2.Bold = 1
$closeQuote2.Bold = 0

Write-Output "Dialogue updated."
